$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "Valentina Perghem "
$ws.Range("B35").Value = "Matteo Zanlucchi | SBARX"
$ws.Range("C35").Value = "Luca Frasca | Clitoriders"
$ws.Range("D35").Value = "Alessio Bragagna | SHARK ATTACK"
$ws.Range("E35").Value = "Michele Leonardi | Rita Levi’s"
$ws.Range("F35").Value = "Matteo Giovannella | Bevem4tut"
